# Form the consolidated report: recompute the "Absent" column (H) for each
# attendance row from the "Real" column (E). Absent = 1 - Real.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 3; $row -le 21; $row++) {
    $real = $ws.Cells.Item($row, 5).Value2   # column E = Real
    $ws.Cells.Item($row, 8).Value = 1 - $real   # column H = Absent
}
